$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (2026/02/25, 水, 5, 48) was inserted at row 853, pushing the
# existing rows 853-894 down to 854-895.
$ws.Rows(853).Insert()

# Column A holds a date-like string ("2026/02/25") that must stay plain text
# (matching the original inlineStr cells) instead of being auto-converted to
# a date serial number by Excel's type inference. Force text format, assign
# the value, then strip the temporary formatting so the cell ends up with no
# explicit style, just like its neighbours.
$ws.Range("A853").NumberFormat = "@"
$ws.Range("A853").Value = "2026/02/25"
$ws.Range("A853").ClearFormats()

$ws.Range("B853").Value = "水"
$ws.Range("C853").Value = 5
$ws.Range("D853").Value = 48
